# Refresh export reports with latest season data
# - Season_Scores: bump Total_Score / Total_Wins for the users who picked up
#   an extra week of results
# - User_Picks: append 4 more weeks' worth of pick rows (weeks 37-56, a
#   repeat of the same 5-user block already present)
# - Summary: bump Total Weeks and refresh the Export Date stamp

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Season_Scores sheet - update Total_Score (G) / Total_Wins (H) for a few users
# ---------------------------------------------------------------------------
$scores = $wb.Worksheets.Item("Season_Scores")

$scores.Cells.Item(4, 7).Value = 99   # G4 Total_Score  (kvstabla)      63 -> 99
$scores.Cells.Item(4, 8).Value = 11   # H4 Total_Wins   (kvstabla)       7 -> 11
$scores.Cells.Item(5, 7).Value = 77   # G5 Total_Score  (RStrype)       49 -> 77
$scores.Cells.Item(8, 7).Value = 77   # G8 Total_Score  (kushminzada)   49 -> 77
$scores.Cells.Item(26, 7).Value = 11  # G26 Total_Score (nalaknas)       7 -> 11
$scores.Cells.Item(33, 7).Value = 88  # G33 Total_Score (PSekhar7)      56 -> 88

# ---------------------------------------------------------------------------
# User_Picks sheet - append rows 37-56 (four more repeats of the 5-row block)
# ---------------------------------------------------------------------------
$picks = $wb.Worksheets.Item("User_Picks")

# Columns B (Season) and C (User_ID) hold numeric-looking text in this sheet
# (e.g. "2025", "1005271556681232384") - format the target range as Text
# first so Excel doesn't silently coerce them into numbers.
$picks.Range("B37:C56").NumberFormat = "@"

$pattern = @(
    @("1005271556681232384", 9),
    @("997251019148951552", 8),
    @("1125900327032598528", 7),
    @("1005301584701521920", 7),
    @("865061949451403264", 1)
)

$row = 37
for ($block = 0; $block -lt 4; $block++) {
    foreach ($entry in $pattern) {
        $picks.Cells.Item($row, 1).Value = 1          # A Week
        $picks.Cells.Item($row, 2).Value = "2025"      # B Season
        $picks.Cells.Item($row, 3).Value = $entry[0]   # C User_ID
        $picks.Cells.Item($row, 4).Value = "N/A"       # D Picks
        $picks.Cells.Item($row, 5).Value = $entry[1]   # E Score
        $row++
    }
}

# Drop the temporary text format back to Normal style so we don't leave the
# new rows with a different cell style than the rest of the sheet.
$picks.Range("B37:C56").Style = "Normal"

# ---------------------------------------------------------------------------
# Summary sheet - refresh Total Weeks + Export Date
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Cells.Item(2, 2).Value = 12                       # B2 Total Weeks   8 -> 12
$summary.Cells.Item(5, 2).Value = "2025-09-09 13:20:27"     # B5 Export Date
